$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.ClearFormats()
}

Set-TextValue 'D2' '26.503.70'
Set-TextValue 'E2' '  +0.70%  '
Set-TextValue 'D3' '1.838.69'
Set-TextValue 'E3' '  +0.15%  '
Set-TextValue 'E4' '  -0.05%  '
Set-TextValue 'D5' '260.00'
Set-TextValue 'D6' '1.000'
Set-TextValue 'E6' '  -0.09%  '
Set-TextValue 'D7' '0.5261'
Set-TextValue 'E7' '  +1.13%  '
Set-TextValue 'D8' '0.3206'
Set-TextValue 'E8' '  -0.59%  '
Set-TextValue 'D9' '0.06783'
Set-TextValue 'E9' '  +0.66%  '
Set-TextValue 'E10' '  +1.57%  '
Set-TextValue 'D11' '0.7832'
Set-TextValue 'E11' '  +3.21%  '
Set-TextValue 'D12' '0.07737'
Set-TextValue 'E12' '  +1.21%  '
Set-TextValue 'D13' '1.842.06'
Set-TextValue 'E13' '  +0.51%  '
Set-TextValue 'D14' '87.64'
Set-TextValue 'E14' '  -0.69%  '
Set-TextValue 'D15' '5.007'
Set-TextValue 'E15' '  -0.02%  '
Set-TextValue 'E16' '  -0.05%  '
Set-TextValue 'D17' '13.82'
Set-TextValue 'E17' '  -0.54%  '
Set-TextValue 'E18' '  -0.01%  '
Set-TextValue 'D19' '0.000007920'
Set-TextValue 'E19' '  +0.68%  '
Set-TextValue 'D20' '26.526.08'
Set-TextValue 'E20' '  +0.66%  '
Set-TextValue 'D21' '2.076.45'
Set-TextValue 'E21' '  +0.26%  '
Set-TextValue 'D22' '4.616'
Set-TextValue 'E22' '  +1.70%  '
Set-TextValue 'E23' '  +1.08%  '
Set-TextValue 'E24' '  -0.35%  '
Set-TextValue 'D25' '141.39'
Set-TextValue 'E25' '  -1.87%  '
Set-TextValue 'D26' '2.148'
Set-TextValue 'E26' '  -3.79%  '
Set-TextValue 'D27' '1.681'
Set-TextValue 'E27' '  +2.25%  '
Set-TextValue 'E28' '  +0.04%  '
Set-TextValue 'D29' '111.82'
Set-TextValue 'E29' '  +0.52%  '
Set-TextValue 'D30' '4.142'
Set-TextValue 'E30' '  -0.43%  '
Set-TextValue 'D31' '0.08692'
Set-TextValue 'E31' '  -0.16%  '
Set-TextValue 'D32' '4.060'
Set-TextValue 'E32' '  -1.66%  '
Set-TextValue 'D33' '0.04858'
Set-TextValue 'E33' '  +1.98%  '
Set-TextValue 'D34' '0.7291'
Set-TextValue 'E34' '  +5.13%  '
Set-TextValue 'D35' '1.132'
Set-TextValue 'E35' '  +1.98%  '
Set-TextValue 'D36' '2.854'
Set-TextValue 'E36' '  +0.01%  '
Set-TextValue 'D37' '3.087'
Set-TextValue 'E37' '  +1.05%  '
Set-TextValue 'D38' '2.240'
Set-TextValue 'E38' '  +2.19%  '
Set-TextValue 'E39' '  +0.31%  '
Set-TextValue 'D40' '0.4775'
Set-TextValue 'E40' '  -0.83%  '
Set-TextValue 'D41' '0.8911'
Set-TextValue 'E41' '  +1.05%  '
Set-TextValue 'D42' '109.64'
Set-TextValue 'E42' '  -1.06%  '
Set-TextValue 'D43' '5.909'
Set-TextValue 'E43' '  -2.97%  '
Set-TextValue 'E44' '  -0.02%  '
Set-TextValue 'D45' '7.643'
Set-TextValue 'E45' '  +0.14%  '
Set-TextValue 'D46' '0.4143'
Set-TextValue 'E46' '  +0.59%  '
Set-TextValue 'B47' 'EnergySwap'
Set-TextValue 'C47' 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue 'D47' '8.975'
Set-TextValue 'E47' '  +0.63%  '
Set-TextValue 'B48' 'Cronos'
Set-TextValue 'C48' 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue 'D48' '0.05850'
Set-TextValue 'E48' '  +0.11%  '
Set-TextValue 'D49' '0.1232'
Set-TextValue 'E49' '  +0.32%  '
Set-TextValue 'D50' '34.81'
Set-TextValue 'E50' '  +0.63%  '
Set-TextValue 'D51' '0.8946'
Set-TextValue 'E51' '  +1.56%  '
